$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

$ws_ALC.Range("H112").Value = 7372250.5
$ws_ALC.Range("J112").Value = 7793510.5
$ws_ALC.Range("L112").Value = 23380531.5
$ws_ALC.Range("N112").Value = -23382747.5

$ws_ALC.Range("H137").Value = 28572986
$ws_ALC.Range("I137").Value = 37038172
$ws_ALC.Range("K137").Value = 111114516
$ws_ALC.Range("M137").Value = -111111966

$ws_ARM.Range("H61").Value = 2650.8635
$ws_ARM.Range("I61").Value = 1370.4
$ws_ARM.Range("J61").Value = 4335.684
$ws_ARM.Range("K61").Value = 1370.4
$ws_ARM.Range("L61").Value = 4335.684
$ws_ARM.Range("M61").Value = -1158.4
$ws_ARM.Range("N61").Value = -4759.684

$ws_ARM.Range("H74").Value = 4723.3057
$ws_ARM.Range("I74").Value = 1402.2222
$ws_ARM.Range("J74").Value = 8044.3887
$ws_ARM.Range("K74").Value = 1402.2222
$ws_ARM.Range("L74").Value = 8044.3887
$ws_ARM.Range("M74").Value = -528.2221999999999
$ws_ARM.Range("N74").Value = -9792.3887

$ws_ARM.Range("H77").Value = 4723.3057
$ws_ARM.Range("I77").Value = 1402.2222
$ws_ARM.Range("J77").Value = 8044.3887
$ws_ARM.Range("K77").Value = 7011.111
$ws_ARM.Range("L77").Value = 40221.9435
$ws_ARM.Range("M77").Value = -2643.111
$ws_ARM.Range("N77").Value = -48957.9435

$ws_ARM.Range("H133").Value = 34666.5
$ws_ARM.Range("J133").Value = 34666.5
$ws_ARM.Range("L133").Value = 34666.5
$ws_ARM.Range("N133").Value = -39726.5

$ws_ARM.Range("H136").Value = 2650.8635
$ws_ARM.Range("I136").Value = 1370.4
$ws_ARM.Range("J136").Value = 4335.684
$ws_ARM.Range("K136").Value = 4111.200000000001
$ws_ARM.Range("L136").Value = 13007.052
$ws_ARM.Range("M136").Value = -1561.200000000001
$ws_ARM.Range("N136").Value = -18107.052

$ws_BSM.Range("H134").Value = 3661.4614
$ws_BSM.Range("I134").Value = 2913.4736
$ws_BSM.Range("J134").Value = 5691.7144
$ws_BSM.Range("K134").Value = 8740.4208
$ws_BSM.Range("L134").Value = 17075.1432
$ws_BSM.Range("M134").Value = -6205.4208
$ws_BSM.Range("N134").Value = -22145.1432

$ws_CRP.Range("H31").Value = 1314.5172
$ws_CRP.Range("I31").Value = 1070.9546
$ws_CRP.Range("J31").Value = 2080
$ws_CRP.Range("K31").Value = 1070.9546
$ws_CRP.Range("L31").Value = 2080
$ws_CRP.Range("M31").Value = -775.9546
$ws_CRP.Range("N31").Value = -2670

$ws_CRP.Range("H34").Value = 1314.5172
$ws_CRP.Range("I34").Value = 1070.9546
$ws_CRP.Range("J34").Value = 2080
$ws_CRP.Range("K34").Value = 1070.9546
$ws_CRP.Range("L34").Value = 2080
$ws_CRP.Range("M34").Value = -868.9546
$ws_CRP.Range("N34").Value = -2484

$ws_CRP.Range("H58").Value = 1874.3256
$ws_CRP.Range("I58").Value = 1252.3478
$ws_CRP.Range("J58").Value = 2589.6
$ws_CRP.Range("K58").Value = 1252.3478
$ws_CRP.Range("L58").Value = 2589.6
$ws_CRP.Range("M58").Value = -1049.3478
$ws_CRP.Range("N58").Value = -2995.6

$ws_CRP.Range("H132").Value = 2115.027
$ws_CRP.Range("I132").Value = 1759.8276
$ws_CRP.Range("J132").Value = 3402.625
$ws_CRP.Range("K132").Value = 5279.4828
$ws_CRP.Range("L132").Value = 10207.875
$ws_CRP.Range("M132").Value = -2749.4828
$ws_CRP.Range("N132").Value = -15267.875

$ws_CRP.Range("H134").Value = 2534.9688
$ws_CRP.Range("I134").Value = 1282.3889
$ws_CRP.Range("J134").Value = 4145.4287
$ws_CRP.Range("K134").Value = 3847.1667
$ws_CRP.Range("L134").Value = 12436.2861
$ws_CRP.Range("M134").Value = -1312.1667
$ws_CRP.Range("N134").Value = -17506.2861

$ws_CRP.Range("H136").Value = 1874.3256
$ws_CRP.Range("I136").Value = 1252.3478
$ws_CRP.Range("J136").Value = 2589.6
$ws_CRP.Range("K136").Value = 3757.0434
$ws_CRP.Range("L136").Value = 7768.799999999999
$ws_CRP.Range("M136").Value = -1207.0434
$ws_CRP.Range("N136").Value = -12868.8

$ws_CRP.Range("H137").Value = 24695
$ws_CRP.Range("J137").Value = 29593.334
$ws_CRP.Range("L137").Value = 29593.334
$ws_CRP.Range("N137").Value = -39793.334

$ws_CRP.Range("H138").Value = 43926.668
$ws_CRP.Range("J138").Value = 43926.668
$ws_CRP.Range("L138").Value = 43926.668
$ws_CRP.Range("N138").Value = -54206.668

$ws_CRP.Range("H140").Value = 37868
$ws_CRP.Range("J140").Value = 44835
$ws_CRP.Range("L140").Value = 44835
$ws_CRP.Range("N140").Value = -55195

$ws_CUL.Range("H113").Value = 8197383
$ws_CUL.Range("I113").Value = 493.76923
$ws_CUL.Range("J113").Value = 10417374
$ws_CUL.Range("K113").Value = 1481.30769
$ws_CUL.Range("L113").Value = 31252122
$ws_CUL.Range("M113").Value = 688.6923099999999
$ws_CUL.Range("N113").Value = -31256462

$ws_GSM.Range("H126").Value = 2462.8928
$ws_GSM.Range("I126").Value = 1844.3846
$ws_GSM.Range("J126").Value = 2998.9333
$ws_GSM.Range("K126").Value = 5533.1538
$ws_GSM.Range("L126").Value = 8996.7999
$ws_GSM.Range("M126").Value = -3063.1538
$ws_GSM.Range("N126").Value = -13936.7999

$ws_GSM.Range("H137").Value = 42755.8
$ws_GSM.Range("J137").Value = 42755.8
$ws_GSM.Range("L137").Value = 42755.8
$ws_GSM.Range("N137").Value = -52955.8

$ws_GSM.Range("H138").Value = 35199.855
$ws_GSM.Range("J138").Value = 35199.855
$ws_GSM.Range("L138").Value = 35199.855
$ws_GSM.Range("N138").Value = -45479.855

$ws_LTW.Range("H132").Value = 4178.8477
$ws_LTW.Range("I132").Value = 3699.1538
$ws_LTW.Range("J132").Value = 4802.45
$ws_LTW.Range("K132").Value = 11097.4614
$ws_LTW.Range("L132").Value = 14407.35
$ws_LTW.Range("M132").Value = -8567.4614
$ws_LTW.Range("N132").Value = -19467.35

$ws_LTW.Range("H133").Value = 46731.5
$ws_LTW.Range("J133").Value = 46731.5
$ws_LTW.Range("L133").Value = 46731.5
$ws_LTW.Range("N133").Value = -51791.5

$ws_LTW.Range("H135").Value = 31428
$ws_LTW.Range("J135").Value = 31428
$ws_LTW.Range("L135").Value = 31428
$ws_LTW.Range("N135").Value = -41568

$ws_LTW.Range("H136").Value = 3260.8572
$ws_LTW.Range("I136").Value = 2354.6829
$ws_LTW.Range("J136").Value = 4949.636
$ws_LTW.Range("K136").Value = 7064.048699999999
$ws_LTW.Range("L136").Value = 14848.908
$ws_LTW.Range("M136").Value = -4514.048699999999
$ws_LTW.Range("N136").Value = -19948.908

$ws_LTW.Range("H141").Value = 42500
$ws_LTW.Range("J141").Value = 42500
$ws_LTW.Range("L141").Value = 42500
$ws_LTW.Range("N141").Value = -52860

$ws_WVR.Range("H132").Value = 45462564
$ws_WVR.Range("I132").Value = 55565080
$ws_WVR.Range("J132").Value = 1255
$ws_WVR.Range("K132").Value = 166695240
$ws_WVR.Range("L132").Value = 3765
$ws_WVR.Range("M132").Value = -166692710
$ws_WVR.Range("N132").Value = -8825

$ws_WVR.Range("H136").Value = 7599813.5
$ws_WVR.Range("I136").Value = 11529660
$ws_WVR.Range("J136").Value = 2109.8
$ws_WVR.Range("K136").Value = 34588980
$ws_WVR.Range("L136").Value = 6329.400000000001
$ws_WVR.Range("M136").Value = -34586430
$ws_WVR.Range("N136").Value = -11429.4

$ws_WVR.Range("H137").Value = 39000
$ws_WVR.Range("J137").Value = 39000
$ws_WVR.Range("L137").Value = 39000
$ws_WVR.Range("N137").Value = -49200

$ws_WVR.Range("H139").Value = 43600
$ws_WVR.Range("J139").Value = 43600
$ws_WVR.Range("L139").Value = 43600
$ws_WVR.Range("N139").Value = -53880

$ws_WVR.Range("H141").Value = 59500
$ws_WVR.Range("J141").Value = 59500
$ws_WVR.Range("L141").Value = 59500
$ws_WVR.Range("N141").Value = -69860
